# Generate Report for Handback
#
# The handback report records, per localized file, the datetime the file
# was handed off and the datetime it was handed back. A new handback run
# completed for the "8e8a11af-...md" file (row 3) in both the zh-cn and
# de-de sheets, so its "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) get refreshed with the newly
# recorded timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 3 (8e8a11af-...md)
$wsZhCn.Range("D3").Value = "2016-01-25 06:31:49"
$wsZhCn.Range("G3").Value = "2016-01-25 06:32:33"

# de-de sheet, row 3 (8e8a11af-...md)
$wsDeDe.Range("D3").Value = "2016-01-25 06:32:00"
$wsDeDe.Range("G3").Value = "2016-01-25 06:32:52"
